# Swap the order of "dnasr281@gmail.com, System" to "System, dnasr281@gmail.com"
# in column G, for every row where that exact text appears.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

$oldValue = "dnasr281@gmail.com, System"
$newValue = "System, dnasr281@gmail.com"

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

for ($row = 1; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, 7)  # Column G is the 7th column
    if ($cell.Value2 -eq $oldValue) {
        $cell.Value2 = $newValue
    }
}
